$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to file 9f02a67f-95fa-46da-bed3-124b9cf5effb.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-17-18 18:17:05"

# --- zh-cn sheet: row 3 corresponds to file 9f02a67f-95fa-46da-bed3-124b9cf5effb.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-18 18:17:01"
$wsZhCn.Range("E3").Value = "9f02a67f-95fa-46da-bed3-124b9cf5effb.c4ab08dd2ba772b5de98b45d661c00d1c6472e85.zh-cn.xlf"

# --- de-de sheet: row 3 corresponds to file 9f02a67f-95fa-46da-bed3-124b9cf5effb.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 18:17:05"
